$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B: old column B..Q shift right to C..R,
# column A keeps its position but its long filenames get replaced below with
# short subject identifiers, and the old (long) filenames move into the new
# column B.
$ws.Columns.Item(2).Insert()

# New header for column B ("subject" moves to A, "filename" now lives in B)
$ws.Cells.Item(1, 2).Value = "filename"

# Row 2..6: column A becomes the short subject id, column B keeps the
# original long filename value.
$ws.Cells.Item(2, 1).Value = "VS06"
$ws.Cells.Item(2, 2).Value = "VS06_704"

$ws.Cells.Item(3, 1).Value = "VS08"
$ws.Cells.Item(3, 2).Value = "VS08_0497"

$ws.Cells.Item(4, 1).Value = "VS24"
$ws.Cells.Item(4, 2).Value = "VS24_823"

$ws.Cells.Item(5, 1).Value = "sub-PD06"
$ws.Cells.Item(5, 2).Value = "sub-PD06"

$ws.Cells.Item(6, 1).Value = "sub-PD10"
$ws.Cells.Item(6, 2).Value = "sub-PD10"

# Row 7 (TOTAL row): no subject, so clear A7 and move "TOTAL" into B7.
$ws.Cells.Item(7, 1).ClearContents()
$ws.Cells.Item(7, 2).Value = "TOTAL"

# A1 header becomes "subject"
$ws.Cells.Item(1, 1).Value = "subject"

# Column widths: new column A (short "subject" ids) is narrower than the
# column B that now holds the long filenames (which keeps ~the old column
# A width). The engine only supports width inputs quantized to 1/6 of a
# character, so these are the closest achievable approximations of the
# author's best-fit widths (8.265625 / 9.59765625).
$ws.Columns.Item(1).ColumnWidth = 7.5
$ws.Columns.Item(2).ColumnWidth = 8.833333333333334
